$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$data = @(
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45219, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 50, 4000, 4000, 4000, '$/paquete 36 unidades', 'Región Metropolitana', 111, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45219, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 3000, 3000, 3000, '$/paquete 36 unidades', 'Región Metropolitana', 83, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44727, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 450, 6000, 6500, 6222, '$/paquete 36 unidades', 'Región Metropolitana', 173, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44757, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44757, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45133, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 4500, 5000, 4750, '$/paquete 36 unidades', 'Región Metropolitana', 132, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44945, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 150, 3500, 3800, 3640, '$/paquete 36 unidades', 'Región Metropolitana', 101, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44904, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 250, 3000, 3200, 3104, '$/paquete 36 unidades', 'Región Metropolitana', 86, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45043, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 90, 3300, 3500, 3389, '$/paquete 36 unidades', 'Región Metropolitana', 94, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44701, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 180, 6000, 6500, 6222, '$/paquete 36 unidades', 'Región Metropolitana', 173, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45007, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45007, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44777, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 150, 7500, 8000, 7767, '$/paquete 36 unidades', 'Región Metropolitana', 216, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44321, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44321, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44944, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 3000, 3200, 3091, '$/paquete 36 unidades', 'Región Metropolitana', 86, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45072, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 110, 5000, 5500, 5273, '$/paquete 36 unidades', 'Región Metropolitana', 146, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44951, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44951, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44188, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44188, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44943, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 450, 2700, 2800, 2756, '$/paquete 36 unidades', 'Región Metropolitana', 77, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45118, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 30, 5500, 5500, 5500, '$/paquete 36 unidades', 'Región Metropolitana', 153, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45090, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 3200, 3500, 3350, '$/paquete 36 unidades', 'Región Metropolitana', 93, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44644, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 160, 6500, 7000, 6750, '$/paquete 36 unidades', 'Región Metropolitana', 188, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44637, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 110, 6500, 7000, 6773, '$/paquete 36 unidades', 'Región Metropolitana', 188, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44679, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44679, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44335, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 150, 600, 700, 633, '$/paquete 6 unidades', 'Región de Ñuble', 106, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44335, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44806, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 4000, 4500, 4227, '$/paquete 36 unidades', 'Región Metropolitana', 117, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45030, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 130, 5000, 6000, 5385, '$/paquete 36 unidades', 'Región Metropolitana', 150, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44328, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44328, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45217, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 80, 4000, 4000, 4000, '$/paquete 36 unidades', 'Región Metropolitana', 111, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44308, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45112, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 80, 7000, 7500, 7188, '$/paquete 36 unidades', 'Región Metropolitana', 200, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44937, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44937, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45106, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 60, 7500, 8000, 7750, '$/paquete 36 unidades', 'Región Metropolitana', 215, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44554, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44554, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44631, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 6000, 6500, 6227, '$/paquete 36 unidades', 'Región Metropolitana', 173, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45049, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 40, 4500, 5000, 4750, '$/paquete 36 unidades', 'Región Metropolitana', 132, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45205, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 70, 4000, 4500, 4214, '$/paquete 36 unidades', 'Región Metropolitana', 117, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44771, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44771, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44839, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44687, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 8000, 8500, 8273, '$/paquete 36 unidades', 'Región Metropolitana', 230, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45071, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 80, 5000, 5000, 5000, '$/paquete 36 unidades', 'Región Metropolitana', 139, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45014, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 5000, 5000, 5000, '$/paquete 36 unidades', 'Región Metropolitana', 139, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45014, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 3500, 3500, 3500, '$/paquete 36 unidades', 'Región Metropolitana', 97, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44616, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44763, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 80, 5000, 5500, 5188, '$/paquete 36 unidades', 'Región Metropolitana', 144, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45036, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 60, 4000, 4500, 4250, '$/paquete 36 unidades', 'Región Metropolitana', 118, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44848, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 2000, 2500, 2250, '$/paquete 36 unidades', 'Región Metropolitana', 62, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44694, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45174, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 50, 4500, 4500, 4500, '$/paquete 36 unidades', 'Región Metropolitana', 125, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44643, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 180, 6500, 7000, 6778, '$/paquete 36 unidades', 'Región Metropolitana', 188, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45149, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 50, 4500, 5000, 4700, '$/paquete 36 unidades', 'Región Metropolitana', 131, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44729, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 6000, 6500, 6273, '$/paquete 36 unidades', 'Región Metropolitana', 174, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44981, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 6000, 6500, 6250, '$/paquete 36 unidades', 'Región Metropolitana', 174, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45076, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 4000, 4200, 4100, '$/paquete 36 unidades', 'Región Metropolitana', 114, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44972, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44972, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44491, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 600, 700, 650, '$/paquete 6 unidades', 'Región Metropolitana', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44491, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 500, 500, 500, '$/paquete 6 unidades', 'Región Metropolitana', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44230, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 600, 700, 650, '$/paquete 6 unidades', 'Región de Ñuble', 108, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44230, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 50, 500, 500, 500, '$/paquete 6 unidades', 'Región de Ñuble', 83, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45028, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 180, 5500, 6000, 5722, '$/paquete 36 unidades', 'Región Metropolitana', 159, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44993, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 300, 600, 800, 700, '$/paquete 6 unidades', 'Región Metropolitana', 117, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44698, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 350, 7500, 8000, 7786, '$/paquete 36 unidades', 'Región Metropolitana', 216, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44797, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 700, 800, 750, '$/paquete 6 unidades', 'Región de Ñuble', 125, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44797, 8, 100112037, 'Cebollín', 'Sin especificar', 'Segunda', 100, 600, 600, 600, '$/paquete 6 unidades', 'Región de Ñuble', 100, 6, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44714, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 180, 6500, 7000, 6778, '$/paquete 36 unidades', 'Región Metropolitana', 188, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44775, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 200, 7000, 7500, 7250, '$/paquete 36 unidades', 'Región Metropolitana', 201, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45063, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 40, 3500, 4000, 3750, '$/paquete 36 unidades', 'Región Metropolitana', 104, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 44649, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 220, 8000, 8500, 8227, '$/paquete 36 unidades', 'Región Metropolitana', 229, 36, 'Hortaliza'),
  @(11, 'Vega Monumental Concepción', 'Bíobío', 45216, 8, 100112037, 'Cebollín', 'Sin especificar', 'Primera', 100, 4000, 4500, 4250, '$/paquete 36 unidades', 'Región Metropolitana', 118, 36, 'Hortaliza')
)

$startRow = 64
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($j = 0; $j -lt $row.Count; $j++) {
    $ws.Cells.Item($r, $j+1).Value = $row[$j]
  }
}

# Match the date-column (D) number format on the two brand-new rows (144, 145)
# to the existing date cells, since they did not exist in the original sheet.
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Range("D144").NumberFormat = $dateFormat
$ws.Range("D145").NumberFormat = $dateFormat